$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (PMC7273137 / elsevier) - author list gains CROSSREF-era spacing, location tag gets _CROSSREF suffix
$ws.Range("E2").Value = "[Ping%Zhang%NULL%1,                    Zhigang%He%NULL%1,                    Gang%Yu%NULL%1,                    Dan%Peng%NULL%1,                    Yikuan%Feng%NULL%1,                    Jianmin%Ling%NULL%1,                    Ye%Wang%NULL%1,                    Shusheng%Li%NULL%0,                    Yi%Bian%NULL%1]"
$ws.Range("I2").Value = "_PMC_elsevier_CROSSREF"

# Row 3 (PMC7175450 / Springer) - same treatment
$ws.Range("E3").Value = "[Tao%Li%NULL%1,                    Yalan%Zhang%NULL%2,                    Yalan%Zhang%NULL%0,                    Cheng%Gong%NULL%1,                    Jing%Wang%NULL%0,                    Bao%Liu%NULL%1,                    Li%Shi%NULL%1,                    Jun%Duan%junjununun@163.com%1]"
$ws.Range("I3").Value = "_PMC_Springer_CROSSREF"

# Row 4 (10.1002/jmv.25796) - now resolved via CROSSREF: title, authors, doi-as-id, format, accepted date
$ws.Range("C4").Value = "`"Caution should be exercised for the detection of SARS\u2010CoV\u20102, especially in the elderly`""
$ws.Range("E4").Value = "[Yajun%Yuan%xref no email%1, Nan%Wang%xref no email%1, Xueqing%Ou%xref no email%1]"
$ws.Range("F4").Value = "10.1002/jmv.25796"
$ws.Range("G4").Value = "CROSSREF"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "2023-01-28"

# Row 5 (PMC7220650 / Springer) - same treatment as rows 2 & 3
$ws.Range("E5").Value = "[Gaoli%Liu%NULL%1,                    Shaowen%Zhang%NULL%2,                    Shaowen%Zhang%NULL%0,                    Zhangfan%Mao%NULL%1,                    Weixing%Wang%13392186@qq.com%1,                    Haifeng%Hu%NULL%1]"
$ws.Range("I5").Value = "_PMC_Springer_CROSSREF"
